$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 3 new rows before row 225. This pushes the existing rows
#    225-244 down to 228-247, carrying their content/format untouched
#    (matches the growth of the sheet dimension from A1:T244 to A1:T247).
$ws.Rows("225:227").Insert()

# 2) The 3 freshly-inserted rows 225-227 are blank; populate them with a
#    new weekly price-report batch (Mercado/Producto columns match the
#    rest of the sheet, dated 2021-09-10 / serial 44449, origin Brasil)
#    for the three usual Calidad tiers: Especial / Primera / Segunda.

# Row 225 - Calidad "Especial"
$ws.Range("A225").Value2 = 9
$ws.Range("B225").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C225").Value2 = "Metropolitana"
$ws.Range("D225").Value2 = 44449
$ws.Range("E225").Value2 = 13
$ws.Range("F225").Value2 = "Fruta"
$ws.Range("G225").Value2 = 100108
$ws.Range("H225").Value2 = "Tropicales y subtropicales"
$ws.Range("I225").Value2 = 100108002
$ws.Range("J225").Value2 = "Mango"
$ws.Range("K225").Value2 = "Sin especificar"
$ws.Range("L225").Value2 = "Especial"
$ws.Range("M225").Value2 = 120
$ws.Range("N225").Value2 = 9000
$ws.Range("O225").Value2 = 9000
$ws.Range("P225").Value2 = 9000
$ws.Range("Q225").Value2 = "$/bandeja 4 kilos"
$ws.Range("R225").Value2 = "Brasil"
$ws.Range("S225").Value2 = 2250
$ws.Range("T225").Value2 = 4

# Row 226 - Calidad "Primera"
$ws.Range("A226").Value2 = 9
$ws.Range("B226").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C226").Value2 = "Metropolitana"
$ws.Range("D226").Value2 = 44449
$ws.Range("E226").Value2 = 13
$ws.Range("F226").Value2 = "Fruta"
$ws.Range("G226").Value2 = 100108
$ws.Range("H226").Value2 = "Tropicales y subtropicales"
$ws.Range("I226").Value2 = 100108002
$ws.Range("J226").Value2 = "Mango"
$ws.Range("K226").Value2 = "Sin especificar"
$ws.Range("L226").Value2 = "Primera"
$ws.Range("M226").Value2 = 240
$ws.Range("N226").Value2 = 9000
$ws.Range("O226").Value2 = 9000
$ws.Range("P226").Value2 = 9000
$ws.Range("Q226").Value2 = "$/bandeja 4 kilos"
$ws.Range("R226").Value2 = "Brasil"
$ws.Range("S226").Value2 = 2250
$ws.Range("T226").Value2 = 4

# Row 227 - Calidad "Segunda"
$ws.Range("A227").Value2 = 9
$ws.Range("B227").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C227").Value2 = "Metropolitana"
$ws.Range("D227").Value2 = 44449
$ws.Range("E227").Value2 = 13
$ws.Range("F227").Value2 = "Fruta"
$ws.Range("G227").Value2 = 100108
$ws.Range("H227").Value2 = "Tropicales y subtropicales"
$ws.Range("I227").Value2 = 100108002
$ws.Range("J227").Value2 = "Mango"
$ws.Range("K227").Value2 = "Sin especificar"
$ws.Range("L227").Value2 = "Segunda"
$ws.Range("M227").Value2 = 260
$ws.Range("N227").Value2 = 9000
$ws.Range("O227").Value2 = 9000
$ws.Range("P227").Value2 = 9000
$ws.Range("Q227").Value2 = "$/bandeja 4 kilos"
$ws.Range("R227").Value2 = "Brasil"
$ws.Range("S227").Value2 = 2250
$ws.Range("T227").Value2 = 4
